$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS sheet ---
$ws1 = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws1.Range("C6").Value = 57361.63420026322
$ws1.Range("C7").Value = 57926.017533596576
$ws1.Range("C8").Value = 52133.415780236915
$ws1.Range("C12").Value = 45966.95085431592
$ws1.Range("C13").Value = 45966.95085431592
$ws1.Range("C14").Value = 33096.95085431592
$ws1.Range("C15").Value = 32367.863563315914
$ws1.Range("C16").Value = 31517.653563315922
$ws1.Range("C20").Value = 562525.4700300111
$ws1.Range("C21").Value = 568060.1798458446
$ws1.Range("C22").Value = 511254.1618612602
$ws1.Range("C26").Value = 450781.798595477
$ws1.Range("C27").Value = 450781.798595477
$ws1.Range("C28").Value = 324570.2130954771
$ws1.Range("C29").Value = 317420.30921319197
$ws1.Range("C30").Value = 309082.597316692

# --- HORIZONTAL TAIL sheet ---
$ws4 = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws4.Range("C9").Value = 470.0
$ws4.Range("D9").Value = -14.280503374065285
$ws4.Range("C10").Value = 478.33333333333326
$ws4.Range("D10").Value = -12.760654143108994
